$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 50995.668
$ws.Range("I21").Value = 50995.668
$ws.Range("K21").Value = 50995.668
$ws.Range("M21").Value = -50527.668

$ws.Range("H23").Value = 50995.668
$ws.Range("I23").Value = 50995.668
$ws.Range("K23").Value = 50995.668
$ws.Range("M23").Value = -50761.668

$ws.Range("H31").Value = 3340.6667
$ws.Range("I31").Value = 4998
$ws.Range("J31").Value = 2512
$ws.Range("K31").Value = 14994
$ws.Range("L31").Value = 7536
$ws.Range("M31").Value = -14764
$ws.Range("N31").Value = -7996

$ws.Range("H38").Value = 911.8
$ws.Range("I38").Value = 911.8
$ws.Range("J38").Value = 0
$ws.Range("K38").Value = 2735.4
$ws.Range("L38").Value = 0
$ws.Range("M38").Value = -2363.4
$ws.Range("N38").ClearContents()

$ws.Range("H43").Value = 417401.4
$ws.Range("I43").Value = 17000
$ws.Range("J43").Value = 517501.75
$ws.Range("K43").Value = 17000
$ws.Range("L43").Value = 517501.75
$ws.Range("M43").Value = -16931
$ws.Range("N43").Value = -517639.75

$ws.Range("H99").Value = 608
$ws.Range("I99").Value = 541.5714
$ws.Range("J99").Value = 840.5
$ws.Range("K99").Value = 1624.7142
$ws.Range("L99").Value = 2521.5
$ws.Range("M99").Value = -126.7142000000001
$ws.Range("N99").Value = -5517.5

$ws.Range("H112").Value = 1726.76
$ws.Range("J112").Value = 1949.421
$ws.Range("L112").Value = 5848.263
$ws.Range("N112").Value = -8064.263

$ws.Range("H137").Value = 4795
$ws.Range("I137").Value = 4795
$ws.Range("J137").Value = 0
$ws.Range("K137").Value = 14385
$ws.Range("L137").Value = 0
$ws.Range("M137").Value = -11835
$ws.Range("N137").ClearContents()

$ws.Range("H138").Value = 2601.56
$ws.Range("J138").Value = 2816.3667
$ws.Range("L138").Value = 8449.1001
$ws.Range("N138").Value = -18729.1001

$ws.Range("H141").Value = 11962.2
$ws.Range("I141").Value = 11962.2
$ws.Range("K141").Value = 35886.60000000001
$ws.Range("M141").Value = -30706.60000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 450
$ws.Range("I4").Value = 450
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 450
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = -334
$ws.Range("N4").ClearContents()

$ws.Range("H32").Value = 8340631
$ws.Range("I32").Value = 9439772
$ws.Range("J32").Value = 18559.715
$ws.Range("K32").Value = 9439772
$ws.Range("L32").Value = 18559.715
$ws.Range("M32").Value = -9439485
$ws.Range("N32").Value = -19133.715

$ws.Range("H61").Value = 11033893
$ws.Range("I61").Value = 8337129.5
$ws.Range("J61").Value = 31259624
$ws.Range("K61").Value = 8337129.5
$ws.Range("L61").Value = 31259624
$ws.Range("M61").Value = -8336917.5
$ws.Range("N61").Value = -31260048

$ws.Range("H97").Value = 1685.0344
$ws.Range("I97").Value = 1384.875
$ws.Range("J97").Value = 2054.4614
$ws.Range("K97").Value = 1384.875
$ws.Range("L97").Value = 2054.4614
$ws.Range("M97").Value = -888.875
$ws.Range("N97").Value = -3046.4614

$ws.Range("H132").Value = 3421.2144
$ws.Range("I132").Value = 1706.0322
$ws.Range("J132").Value = 8254.909
$ws.Range("K132").Value = 5118.096600000001
$ws.Range("L132").Value = 24764.727
$ws.Range("M132").Value = -2588.096600000001
$ws.Range("N132").Value = -29824.727

$ws.Range("H136").Value = 11033893
$ws.Range("I136").Value = 8337129.5
$ws.Range("J136").Value = 31259624
$ws.Range("K136").Value = 25011388.5
$ws.Range("L136").Value = 93778872
$ws.Range("M136").Value = -25008838.5
$ws.Range("N136").Value = -93783972

$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 842019.0600000001
$ws.Range("I31").Value = 9456
$ws.Range("J31").Value = 2599652.2
$ws.Range("K31").Value = 9456
$ws.Range("L31").Value = 2599652.2
$ws.Range("M31").Value = -9161
$ws.Range("N31").Value = -2600242.2

$ws.Range("H34").Value = 842019.0600000001
$ws.Range("I34").Value = 9456
$ws.Range("J34").Value = 2599652.2
$ws.Range("K34").Value = 9456
$ws.Range("L34").Value = 2599652.2
$ws.Range("M34").Value = -9254
$ws.Range("N34").Value = -2600056.2

$ws.Range("H58").Value = 1497.5156
$ws.Range("I58").Value = 862.0769
$ws.Range("J58").Value = 4251.0835
$ws.Range("K58").Value = 862.0769
$ws.Range("L58").Value = 4251.0835
$ws.Range("M58").Value = -659.0769
$ws.Range("N58").Value = -4657.0835

$ws.Range("H136").Value = 1497.5156
$ws.Range("I136").Value = 862.0769
$ws.Range("J136").Value = 4251.0835
$ws.Range("K136").Value = 2586.2307
$ws.Range("L136").Value = 12753.2505
$ws.Range("M136").Value = -36.23070000000007
$ws.Range("N136").Value = -17853.2505

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 437.875
$ws.Range("I14").Value = 437.875
$ws.Range("K14").Value = 1313.625
$ws.Range("M14").Value = -1140.625

$ws.Range("H113").Value = 1132
$ws.Range("I113").Value = 491.33334
$ws.Range("J113").Value = 1406.5714
$ws.Range("K113").Value = 1474.00002
$ws.Range("L113").Value = 4219.7142
$ws.Range("M113").Value = 695.9999800000001
$ws.Range("N113").Value = -8559.7142

$ws.Range("H131").Value = 4292.2
$ws.Range("I131").Value = 639.3333
$ws.Range("J131").Value = 5620.515
$ws.Range("K131").Value = 1917.9999
$ws.Range("L131").Value = 16861.545
$ws.Range("M131").Value = 3122.0001
$ws.Range("N131").Value = -26941.545

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 41162.5
$ws.Range("I7").Value = 4686.85
$ws.Range("J7").Value = 132351.62
$ws.Range("K7").Value = 4686.85
$ws.Range("L7").Value = 132351.62
$ws.Range("M7").Value = -4574.85
$ws.Range("N7").Value = -132575.62

$ws.Range("H40").Value = 2194.4443
$ws.Range("I40").Value = 1463.8276
$ws.Range("K40").Value = 1463.8276
$ws.Range("M40").Value = -1327.8276

$ws.Range("H126").Value = 41162.5
$ws.Range("I126").Value = 4686.85
$ws.Range("J126").Value = 132351.62
$ws.Range("K126").Value = 14060.55
$ws.Range("L126").Value = 397054.86
$ws.Range("M126").Value = -11590.55
$ws.Range("N126").Value = -401994.86

$ws.Range("H132").Value = 464756.9
$ws.Range("I132").Value = 12320.765
$ws.Range("J132").Value = 2003039.8
$ws.Range("K132").Value = 36962.295
$ws.Range("L132").Value = 6009119.4
$ws.Range("M132").Value = -34432.295
$ws.Range("N132").Value = -6014179.4

$ws.Range("H136").Value = 102783.62
$ws.Range("I136").Value = 18598.666
$ws.Range("J136").Value = 174942.14
$ws.Range("K136").Value = 55795.99800000001
$ws.Range("L136").Value = 524826.42
$ws.Range("M136").Value = -53245.99800000001
$ws.Range("N136").Value = -529926.42

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 24499.25
$ws.Range("I81").Value = 3500
$ws.Range("J81").Value = 45498.5
$ws.Range("K81").Value = 7000
$ws.Range("L81").Value = 90997
$ws.Range("M81").Value = -5939
$ws.Range("N81").Value = -93119

$ws.Range("H84").Value = 24499.25
$ws.Range("I84").Value = 3500
$ws.Range("J84").Value = 45498.5
$ws.Range("K84").Value = 35000
$ws.Range("L84").Value = 454985
$ws.Range("M84").Value = -29696
$ws.Range("N84").Value = -465593

$ws.Range("H96").Value = 4044.875
$ws.Range("I96").Value = 4044.875
$ws.Range("J96").Value = 0
$ws.Range("K96").Value = 4044.875
$ws.Range("L96").Value = 0
$ws.Range("M96").Value = -2671.875
$ws.Range("N96").ClearContents()

$ws.Range("H107").Value = 41668304
$ws.Range("I107").Value = 62501776
$ws.Range("J107").Value = 1366.5
$ws.Range("K107").Value = 187505328
$ws.Range("L107").Value = 4099.5
$ws.Range("M107").Value = -187503408
$ws.Range("N107").Value = -7939.5

$ws.Range("H132").Value = 3796.9688
$ws.Range("I132").Value = 2836.2
$ws.Range("J132").Value = 7228.2856
$ws.Range("K132").Value = 8508.599999999999
$ws.Range("L132").Value = 21684.8568
$ws.Range("M132").Value = -5978.599999999999
$ws.Range("N132").Value = -26744.8568

$ws.Range("H141").Value = 59500
$ws.Range("J141").Value = 59500
$ws.Range("L141").Value = 59500
$ws.Range("N141").Value = -69860
